# "Execute All Module except ESS"
# Update the employee roster's "Employee" sheet (column A) with the latest
# role-assignment snapshot: a new employee entry is inserted at row 2 and the
# remaining rows are refreshed to their current roster values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

$ws.Range("A2").Value  = "13292186 - Lamont Champlin`nROLE : QCO 1706187312887"
$ws.Range("A4").Value  = "92457737 - Lilliana Williamson`nROLE : RTGO100 1701844270281"
$ws.Range("A5").Value  = "90317880 - Lewis Mosciski`nROLE : RTGO100 1701844270281"
$ws.Range("A6").Value  = "90833312 - Angelo Mueller`nROLE : RTGO100 1701844270281"
$ws.Range("A7").Value  = "92970163 - Glenna Lynch`nROLE : RTGO100 1701853905917"
$ws.Range("A8").Value  = "ROLE GROUP : RTGO Operator 2024-01-09T10:07:14.855048200"
$ws.Range("A9").Value  = "OFF"
